$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "26.149.02"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.677.49"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "213.78"
$ws.Range("E5").Value = "  -1.95%  "
Set-TextValue "D6" "0.5274"
$ws.Range("E6").Value = "  -4.73%  "
Set-TextValue "D7" "1.006"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.88%  "
Set-TextValue "D9" "0.06339"
$ws.Range("E9").Value = "  -2.11%  "
Set-TextValue "D10" "21.31"
$ws.Range("E10").Value = "  -3.29%  "
Set-TextValue "D11" "0.07629"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "1.705.10"
$ws.Range("E12").Value = "  +1.49%  "
Set-TextValue "D13" "4.510"
$ws.Range("E13").Value = "  -0.60%  "
Set-TextValue "D14" "0.5702"
$ws.Range("E14").Value = "  -1.75%  "
Set-TextValue "D15" "0.000008207"
$ws.Range("E15").Value = "  -2.63%  "
Set-TextValue "D16" "65.97"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "26.171.34"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  -1.48%  "
Set-TextValue "D20" "10.66"
$ws.Range("E20").Value = "  -2.12%  "
Set-TextValue "D21" "189.54"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  -0.06%  "
Set-TextValue "D24" "149.11"
$ws.Range("E24").Value = "  +1.64%  "
Set-TextValue "D25" "0.1256"
$ws.Range("E25").Value = "  -4.53%  "
Set-TextValue "D26" "7.655"
$ws.Range("E26").Value = "  -3.05%  "
Set-TextValue "D27" "15.83"
$ws.Range("E27").Value = "  +0.28%  "
Set-TextValue "D28" "0.06401"
$ws.Range("E28").Value = "  +1.10%  "
Set-TextValue "D29" "1.359"
$ws.Range("E29").Value = "  -2.21%  "
Set-TextValue "D30" "1.302"
$ws.Range("E30").Value = "  -1.50%  "
Set-TextValue "D31" "3.547"
$ws.Range("E31").Value = "  -1.09%  "
Set-TextValue "D32" "3.537"
$ws.Range("E32").Value = "  -1.01%  "
Set-TextValue "D33" "1.667"
$ws.Range("E33").Value = "  +0.15%  "
Set-TextValue "D34" "1.015"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "2.418"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.6043"
$ws.Range("E36").Value = "  -2.20%  "
Set-TextValue "D37" "2.716"
$ws.Range("E37").Value = "  +0.02%  "
Set-TextValue "D38" "0.01630"
$ws.Range("E38").Value = "  +0.40%  "
Set-TextValue "D39" "6.139"
$ws.Range("D40").Value = "1.089.47"
$ws.Range("E40").Value = "  -2.11%  "
Set-TextValue "D41" "0.8731"
$ws.Range("E41").Value = "  +0.29%  "
Set-TextValue "D42" "1.009"
$ws.Range("E42").Value = "  -0.57%  "
Set-TextValue "D43" "100.06"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "1.829.84"
$ws.Range("E44").Value = "  -0.02%  "
Set-TextValue "D45" "0.00000000111"
$ws.Range("E45").Value = "  +1.28%  "
Set-TextValue "D46" "57.15"
$ws.Range("E46").Value = "  -0.25%  "
Set-TextValue "D47" "1.004"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.05257"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "8.005"
$ws.Range("E49").Value = "  -2.05%  "
Set-TextValue "D50" "0.4277"
$ws.Range("E50").Value = "  -0.31%  "
Set-TextValue "D51" "5.957"
$ws.Range("E51").Value = "  -1.65%  "
